$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamp values for column A (rows 2-122), taken from the updated
# simulation run (adult#001_clear) - same relative time-step pattern,
# shifted to a new simulation start date/time.
$newTimes = @{
    2 = 44937.5003107185
    3 = 44937.50239405183
    4 = 44937.50447738516
    5 = 44937.5065607185
    6 = 44937.50864405183
    7 = 44937.51072738516
    8 = 44937.51281071849
    9 = 44937.51489405183
    10 = 44937.51697738516
    11 = 44937.51906071849
    12 = 44937.52114405183
    13 = 44937.52322738517
    14 = 44937.5253107185
    15 = 44937.52739405183
    16 = 44937.52947738516
    17 = 44937.5315607185
    18 = 44937.53364405183
    19 = 44937.53572738516
    20 = 44937.5378107185
    21 = 44937.53989405183
    22 = 44937.54197738516
    23 = 44937.54406071849
    24 = 44937.54614405183
    25 = 44937.54822738516
    26 = 44937.55031071849
    27 = 44937.55239405183
    28 = 44937.55447738517
    29 = 44937.5565607185
    30 = 44937.55864405183
    31 = 44937.56072738516
    32 = 44937.5628107185
    33 = 44937.56489405183
    34 = 44937.56697738516
    35 = 44937.5690607185
    36 = 44937.57114405183
    37 = 44937.57322738516
    38 = 44937.57531071849
    39 = 44937.57739405183
    40 = 44937.57947738516
    41 = 44937.58156071849
    42 = 44937.58364405183
    43 = 44937.58572738517
    44 = 44937.5878107185
    45 = 44937.58989405183
    46 = 44937.59197738516
    47 = 44937.5940607185
    48 = 44937.59614405183
    49 = 44937.59822738516
    50 = 44937.6003107185
    51 = 44937.60239405183
    52 = 44937.60447738516
    53 = 44937.60656071849
    54 = 44937.60864405183
    55 = 44937.61072738516
    56 = 44937.61281071849
    57 = 44937.61489405183
    58 = 44937.61697738517
    59 = 44937.6190607185
    60 = 44937.62114405183
    61 = 44937.62322738516
    62 = 44937.6253107185
    63 = 44937.62739405183
    64 = 44937.62947738516
    65 = 44937.6315607185
    66 = 44937.63364405183
    67 = 44937.63572738516
    68 = 44937.63781071849
    69 = 44937.63989405183
    70 = 44937.64197738516
    71 = 44937.64406071849
    72 = 44937.64614405183
    73 = 44937.64822738517
    74 = 44937.6503107185
    75 = 44937.65239405183
    76 = 44937.65447738516
    77 = 44937.6565607185
    78 = 44937.65864405183
    79 = 44937.66072738516
    80 = 44937.6628107185
    81 = 44937.66489405183
    82 = 44937.66697738516
    83 = 44937.66906071849
    84 = 44937.67114405183
    85 = 44937.67322738516
    86 = 44937.67531071849
    87 = 44937.67739405183
    88 = 44937.67947738517
    89 = 44937.6815607185
    90 = 44937.68364405183
    91 = 44937.68572738516
    92 = 44937.6878107185
    93 = 44937.68989405183
    94 = 44937.69197738516
    95 = 44937.6940607185
    96 = 44937.69614405183
    97 = 44937.69822738516
    98 = 44937.70031071849
    99 = 44937.70239405183
    100 = 44937.70447738516
    101 = 44937.70656071849
    102 = 44937.70864405183
    103 = 44937.71072738517
    104 = 44937.7128107185
    105 = 44937.71489405183
    106 = 44937.71697738516
    107 = 44937.7190607185
    108 = 44937.72114405183
    109 = 44937.72322738516
    110 = 44937.7253107185
    111 = 44937.72739405183
    112 = 44937.72947738516
    113 = 44937.73156071849
    114 = 44937.73364405183
    115 = 44937.73572738516
    116 = 44937.73781071849
    117 = 44937.73989405183
    118 = 44937.74197738517
    119 = 44937.7440607185
    120 = 44937.74614405183
    121 = 44937.74822738516
    122 = 44937.7503107185
}

foreach ($row in $newTimes.Keys) {
    $ws.Cells.Item($row, 1).Value = $newTimes[$row]
}
